# Applies the cryptos.xlsx price/volume/listing update described in the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.087.41"
$ws.Range("E2").Value = "  -0.73%  "

# Row 3
$ws.Range("D3").Value = "2.434.55"
$ws.Range("E3").Value = "  -0.29%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.88"
$ws.Range("E5").Value = "  +1.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.95"
$ws.Range("E6").Value = "  -0.63%  "

# Row 8
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("D9").Value = "2.432.92"
$ws.Range("E9").Value = "  -0.19%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  -1.50%  "

# Row 11
$ws.Range("E11").Value = "  +2.14%  "

# Row 12
$ws.Range("E12").Value = "  -0.78%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("E13").Value = "  -2.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.33"
$ws.Range("E14").Value = "  -2.08%  "

# Row 15
$ws.Range("E15").Value = "  -1.53%  "

# Row 16
$ws.Range("D16").Value = "2.877.21"
$ws.Range("E16").Value = "  -0.10%  "

# Row 17
$ws.Range("D17").Value = "61.914.39"
$ws.Range("E17").Value = "  -0.66%  "

# Row 18
$ws.Range("D18").Value = "2.440.44"
$ws.Range("E18").Value = "  +0.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.82"
$ws.Range("E19").Value = "  -3.74%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.17"
$ws.Range("E20").Value = "  -1.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.68"
$ws.Range("E21").Value = "  +0.72%  "

# Row 22
$ws.Range("E22").Value = "  -1.38%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.96"
$ws.Range("E23").Value = "  -4.84%  "

# Row 24
$ws.Range("E24").Value = "  +0.23%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.66"
$ws.Range("E25").Value = "  +0.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.27"
$ws.Range("E26").Value = "  +3.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "606.52"
$ws.Range("E27").Value = "  -1.68%  "

# Row 28
$ws.Range("D28").Value = "2.550.85"
$ws.Range("E28").Value = "  +0.33%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0951"
$ws.Range("E29").Value = "  -3.71%  "

# Row 30
$ws.Range("E30").Value = "  -0.31%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -4.10%  "

# Row 32
$ws.Range("E32").Value = "  -1.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.89"
$ws.Range("E33").Value = "  +0.67%  "

# Row 34
$ws.Range("E34").Value = "  +1.35%  "

# Row 35
$ws.Range("E35").Value = "  -4.12%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.44"
$ws.Range("E37").Value = "  -3.60%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.376"
$ws.Range("E38").Value = "  -0.31%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.34"
$ws.Range("E39").Value = "  +0.22%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.18"
$ws.Range("E40").Value = "  +2.79%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.40"
$ws.Range("E41").Value = "  -1.74%  "

# Row 42
$ws.Range("E42").Value = "  -2.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.68"
$ws.Range("E43").Value = "  +2.22%  "

# Row 44
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("E45").Value = "  -3.87%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.52"
$ws.Range("E46").Value = "  -2.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.62"
$ws.Range("E47").Value = "  -3.22%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.603"
$ws.Range("E48").Value = "  +1.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0522"
$ws.Range("E49").Value = "  -0.98%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0248"
$ws.Range("E50").Value = "  +16.58%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.54"
$ws.Range("E51").Value = "  -5.16%  "
